$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows whose Target cluster is "Inflammatory-Mac"
# (original rows 8 and 9); delete bottom-up so row indices stay valid.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()

# Recomputed TPM-derived values for the remaining rows (2-7).
# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Gdnf"
$ws.Range("C2").Value = "Gfra1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.13541
$ws.Range("H2").Value = 0.40623
$ws.Range("I2").Value = 0.03919062573893041
$ws.Range("J2").Value = 0.0391906257389304
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1389376666666667
$ws.Range("N2").Value = 0.416813
$ws.Range("O2").Value = 0.01722256533596611
$ws.Range("P2").Value = 0.01722256533596611
$ws.Range("Q2").Value = 0.01881354944333333
$ws.Range("R2").Value = 0.16932194499
$ws.Range("S2").Value = 0.000674963112346124
$ws.Range("T2").Value = 0.000674963112346124

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Gdnf"
$ws.Range("C3").Value = "Gfra1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.13541
$ws.Range("H3").Value = 0.40623
$ws.Range("I3").Value = 0.03919062573893041
$ws.Range("J3").Value = 0.0391906257389304
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 4.220039333333333
$ws.Range("N3").Value = 12.660118
$ws.Range("O3").Value = 0.5231115858095611
$ws.Range("P3").Value = 0.5231115858095611
$ws.Range("Q3").Value = 0.5714355261266667
$ws.Range("R3").Value = 5.14291973514
$ws.Range("S3").Value = 0.02050107037916089
$ws.Range("T3").Value = 0.02050107037916088

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.13541
$ws.Range("H4").Value = 0.40623
$ws.Range("I4").Value = 0.03919062573893041
$ws.Range("J4").Value = 0.0391906257389304
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.708210666666667
$ws.Range("N4").Value = 11.124632
$ws.Range("O4").Value = 0.4596658488544727
$ws.Range("P4").Value = 0.4596658488544728
$ws.Range("Q4").Value = 0.5021288063733333
$ws.Range("R4").Value = 4.51915925736
$ws.Range("S4").Value = 0.01801459224742339
$ws.Range("T4").Value = 0.01801459224742339

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.319753
$ws.Range("H5").Value = 9.959258999999999
$ws.Range("I5").Value = 0.9608093742610696
$ws.Range("J5").Value = 0.9608093742610695
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1389376666666667
$ws.Range("N5").Value = 0.416813
$ws.Range("O5").Value = 0.01722256533596611
$ws.Range("P5").Value = 0.01722256533596611
$ws.Range("Q5").Value = 0.4612387357296666
$ws.Range("R5").Value = 4.151148621567
$ws.Range("S5").Value = 0.01654760222361999
$ws.Range("T5").Value = 0.01654760222361999

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Gdnf"
$ws.Range("C6").Value = "Gfra1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.319753
$ws.Range("H6").Value = 9.959258999999999
$ws.Range("I6").Value = 0.9608093742610696
$ws.Range("J6").Value = 0.9608093742610695
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 4.220039333333333
$ws.Range("N6").Value = 12.660118
$ws.Range("O6").Value = 0.5231115858095611
$ws.Range("P6").Value = 0.5231115858095611
$ws.Range("Q6").Value = 14.00948823695133
$ws.Range("R6").Value = 126.085394132562
$ws.Range("S6").Value = 0.5026105154304001
$ws.Range("T6").Value = 0.5026105154304001

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Gdnf"
$ws.Range("C7").Value = "Gfra1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.319753
$ws.Range("H7").Value = 9.959258999999999
$ws.Range("I7").Value = 0.9608093742610696
$ws.Range("J7").Value = 0.9608093742610695
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.708210666666667
$ws.Range("N7").Value = 11.124632
$ws.Range("O7").Value = 0.4596658488544727
$ws.Range("P7").Value = 0.4596658488544728
$ws.Range("Q7").Value = 12.31034348529867
$ws.Range("R7").Value = 110.793091367688
$ws.Range("S7").Value = 0.4416512566070493
$ws.Range("T7").Value = 0.4416512566070493
